$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-01 Thursday" "2024-08-02 Friday"

Replace-Text "63×65=" "61×53="
Replace-Text "85×25=" "83×50="
Replace-Text "98×28=" "84×43="
Replace-Text "86×44=" "94×73="
Replace-Text "25×13=" "71×72="
Replace-Text "77×95=" "68×49="
Replace-Text "35×84=" "60×81="
Replace-Text "49×65=" "62×22="
Replace-Text "78×35=" "65×86="
Replace-Text "26×77=" "94×43="
Replace-Text "84×53=" "92×68="
Replace-Text "56×36=" "54×32="
Replace-Text "65×60=" "27×87="
Replace-Text "63×77=" "39×44="
Replace-Text "77×75=" "96×66="
Replace-Text "17×64=" "55×22="
Replace-Text "87×85=" "16×33="
Replace-Text "35×43=" "75×31="
Replace-Text "71×79=" "74×89="
Replace-Text "61×54=" "83×69="
Replace-Text "65×54=" "56×34="
Replace-Text "98×40=" "23×64="
Replace-Text "56×51=" "79×95="
Replace-Text "82×31=" "45×17="
Replace-Text "42×35=" "48×19="
